# Add the new "Silver Ring of The Labyrinth" raid item as row 9 of Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A9 stays blank (no id yet for this item) but is still a "touched" cell in
# the sheet, same as in the source file. Touching the number format forces
# the cell to exist without giving it a value.
$ws.Range("A9").NumberFormat = "General"

# name / type / specialty_type / description
$ws.Range("C9").Value = "Silver Ring of The Labyrinth"
$ws.Range("D9").Value = "ring"
$ws.Range("F9").Value = "Labyrinth Cloth"
$ws.Range("G9").Value = "A silver ring that glisten with the rage of the little girl"

# base_damage
$ws.Range("I9").Value = 1500

# cost / gold_dust_cost / shards_cost / copper_coin_cost
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 0

# base_damage_mod / base_healing_mod / base_ac_mod
$ws.Range("Q9").Value = 0.33
$ws.Range("R9").Value = 0.33
$ws.Range("S9").Value = 0.33

# agi_mod / focus_mod
$ws.Range("Y9").Value = 0
$ws.Range("Z9").Value = 0

# can_drop
$ws.Range("AC9").Value = 1

# base_damage_mod_bonus / base_healing_mod_bonus / base_ac_mod_bonus /
# fight_time_out_mod_bonus / move_time_out_mod_bonus
$ws.Range("AI9").Value = 0
$ws.Range("AJ9").Value = 0
$ws.Range("AK9").Value = 0
$ws.Range("AL9").Value = 0
$ws.Range("AM9").Value = 0

# kingdom_damage
$ws.Range("AS9").Value = 0

# increase_stat_by
$ws.Range("AV9").Value = 0

# increase_skill_bonus_by / increase_skill_training_bonus_by
$ws.Range("AX9").Value = 0
$ws.Range("AY9").Value = 0

# resurrection_chance
$ws.Range("BA9").Value = 0

# spell_evasion / artifact_annulment / healing_reduction / affix_damage_reduction
$ws.Range("BB9").Value = 1
$ws.Range("BC9").Value = 1
$ws.Range("BD9").Value = 1
$ws.Range("BE9").Value = 1

# devouring_light / devouring_darkness
$ws.Range("BF9").Value = 0
$ws.Range("BG9").Value = 0

# holy_stacks
$ws.Range("BM9").Value = 20

# ambush_chance / ambush_resistance / counter_chance / counter_resistance
$ws.Range("BN9").Value = 0
$ws.Range("BO9").Value = 0
$ws.Range("BP9").Value = 0
$ws.Range("BQ9").Value = 0

# Leave the selection where the author left it when finishing the edit.
$ws.Range("S9").Select()
